$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 6, pushing every
# subsequent record (old rows 6-97) down by one row.
$ws.Rows("6:6").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Vega Monumental Concepción"
$ws.Range("C6").Value = "Bíobío"
$ws.Range("D6").Value = 44552
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100102
$ws.Range("H6").Value = "Cítricos"
$ws.Range("I6").Value = 100102004
$ws.Range("J6").Value = "Mandarina"
$ws.Range("K6").Value = "Murcott"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 310
$ws.Range("N6").Value = 800
$ws.Range("O6").Value = 8500
$ws.Range("P6").Value = 4526
$ws.Range("Q6").Value = "`$/caja 15 kilos"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 302
$ws.Range("T6").Value = 15
